# Update "想去人数" (want-to-go count) values in column F for rows 4-6
# on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1
    $ws.Range("F5").Value = 31
    $ws.Range("F6").Value = 140
}
